# Auto-generated edit script: update Leve profit-calculation values per scheduled market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: H17,J17,L17,N17
$ws.Range("H17").Value = 962.0833
$ws.Range("J17").Value = 962.0833
$ws.Range("L17").Value = 2886.2499
$ws.Range("N17").Value = -3222.2499
# Row 129: H129,I129,J129,K129,L129,M129,N129
$ws.Range("H129").Value = 918.6896400000001
$ws.Range("I129").Value = 310.83334
$ws.Range("J129").Value = 1077.2609
$ws.Range("K129").Value = 932.5000200000001
$ws.Range("L129").Value = 3231.7827
$ws.Range("M129").Value = 4067.49998
$ws.Range("N129").Value = -13231.7827
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 701962.5600000001
$ws.Range("I132").Value = 1915.5927
$ws.Range("J132").Value = 3064621
$ws.Range("K132").Value = 5746.7781
$ws.Range("L132").Value = 9193863
$ws.Range("M132").Value = -3216.7781
$ws.Range("N132").Value = -9198923
# Row 137: H137,I137,K137,M137
$ws.Range("H137").Value = 1725538
$ws.Range("I137").Value = 2326556.2
$ws.Range("K137").Value = 6979668.600000001
$ws.Range("M137").Value = -6977118.600000001
# Row 138: H138,I138,J138,K138,L138,M138,N138
$ws.Range("H138").Value = 2383013.8
$ws.Range("I138").Value = 1543.7073
$ws.Range("J138").Value = 5749920
$ws.Range("K138").Value = 4631.1219
$ws.Range("L138").Value = 17249760
$ws.Range("M138").Value = 508.8780999999999
$ws.Range("N138").Value = -17260040

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32,I32,J32,K32,L32,M32,N32
$ws.Range("H32").Value = 16780.21
$ws.Range("I32").Value = 16847.688
$ws.Range("J32").Value = 16554.305
$ws.Range("K32").Value = 16847.688
$ws.Range("L32").Value = 16554.305
$ws.Range("M32").Value = -16560.688
$ws.Range("N32").Value = -17128.305
# Row 45: H45,I45,K45,M45
$ws.Range("H45").Value = 1404.2858
$ws.Range("I45").Value = 1155
$ws.Range("K45").Value = 1155
$ws.Range("M45").Value = -778
# Row 61: H61,I61,J61,K61,L61,M61,N61
$ws.Range("H61").Value = 27835042
$ws.Range("I61").Value = 40041390
$ws.Range("J61").Value = 93339.27
$ws.Range("K61").Value = 40041390
$ws.Range("L61").Value = 93339.27
$ws.Range("M61").Value = -40041178
$ws.Range("N61").Value = -93763.27
# Row 115: H115,J115,L115,N115
$ws.Range("H115").Value = 29971.273
$ws.Range("J115").Value = 29971.273
$ws.Range("L115").Value = 29971.273
$ws.Range("N115").Value = -33105.273
# Row 122: H122,I122,J122,K122,L122,M122,N122
$ws.Range("H122").Value = 2527060
$ws.Range("I122").Value = 1865.9412
$ws.Range("J122").Value = 11112720
$ws.Range("K122").Value = 5597.8236
$ws.Range("L122").Value = 33338160
$ws.Range("M122").Value = -3147.8236
$ws.Range("N122").Value = -33343060
# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 27835042
$ws.Range("I136").Value = 40041390
$ws.Range("J136").Value = 93339.27
$ws.Range("K136").Value = 120124170
$ws.Range("L136").Value = 280017.81
$ws.Range("M136").Value = -120121620
$ws.Range("N136").Value = -285117.81

$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31,I31,K31,M31
$ws.Range("H31").Value = 2912.3333
$ws.Range("I31").Value = 1652.4
$ws.Range("K31").Value = 1652.4
$ws.Range("M31").Value = -1357.4
# Row 34: H34,I34,K34,M34
$ws.Range("H34").Value = 2912.3333
$ws.Range("I34").Value = 1652.4
$ws.Range("K34").Value = 1652.4
$ws.Range("M34").Value = -1450.4
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 15283.542
$ws.Range("I132").Value = 1187.8214
$ws.Range("K132").Value = 3563.4642
$ws.Range("M132").Value = -1033.4642

$ws = $wb.Worksheets.Item("CUL")
# Row 107: H107,J107,L107,N107
$ws.Range("H107").Value = 822.7143
$ws.Range("J107").Value = 856
$ws.Range("L107").Value = 2568
$ws.Range("N107").Value = -6408

$ws = $wb.Worksheets.Item("GSM")
# Row 70: H70,I70,J70,K70,L70,M70,N70
$ws.Range("H70").Value = 35798.406
$ws.Range("I70").Value = 49647.727
$ws.Range("J70").Value = 5329.9
$ws.Range("K70").Value = 49647.727
$ws.Range("L70").Value = 5329.9
$ws.Range("M70").Value = -49377.727
$ws.Range("N70").Value = -5869.9
# Row 73: H73,I73,J73,K73,L73,M73,N73
$ws.Range("H73").Value = 35798.406
$ws.Range("I73").Value = 49647.727
$ws.Range("J73").Value = 5329.9
$ws.Range("K73").Value = 49647.727
$ws.Range("L73").Value = 5329.9
$ws.Range("M73").Value = -48711.727
$ws.Range("N73").Value = -7201.9
# Row 80: H80,I80,J80,K80,L80,M80,N80
$ws.Range("H80").Value = 3240.25
$ws.Range("I80").Value = 2805.5
$ws.Range("J80").Value = 3481.7778
$ws.Range("K80").Value = 2805.5
$ws.Range("L80").Value = 3481.7778
$ws.Range("M80").Value = -1807.5
$ws.Range("N80").Value = -5477.7778
# Row 83: H83,I83,J83,K83,L83,M83,N83
$ws.Range("H83").Value = 3240.25
$ws.Range("I83").Value = 2805.5
$ws.Range("J83").Value = 3481.7778
$ws.Range("K83").Value = 14027.5
$ws.Range("L83").Value = 17408.889
$ws.Range("M83").Value = -9035.5
$ws.Range("N83").Value = -27392.889
# Row 102: H102,I102,J102,K102,L102,M102,N102
$ws.Range("H102").Value = 1187.2632
$ws.Range("I102").Value = 1072.6666
$ws.Range("J102").Value = 1617
$ws.Range("K102").Value = 1072.6666
$ws.Range("L102").Value = 1617
$ws.Range("M102").Value = 549.3334
$ws.Range("N102").Value = -4861
# Row 107: H107,I107,J107,K107,L107,M107,N107
$ws.Range("H107").Value = 284.44446
$ws.Range("I107").Value = 125.71429
$ws.Range("J107").Value = 840
$ws.Range("K107").Value = 125.71429
$ws.Range("L107").Value = 840
$ws.Range("M107").Value = 1794.28571
$ws.Range("N107").Value = -4680
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 53014.08
$ws.Range("I132").Value = 32735.531
$ws.Range("K132").Value = 98206.59299999999
$ws.Range("M132").Value = -95676.59299999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7: H7,I7,J7,K7,L7,M7,N7
$ws.Range("H7").Value = 5300.222
$ws.Range("I7").Value = 3551
$ws.Range("J7").Value = 6699.6
$ws.Range("K7").Value = 3551
$ws.Range("L7").Value = 6699.6
$ws.Range("M7").Value = -3439
$ws.Range("N7").Value = -6923.6
# Row 40: H40,I40,J40,K40,L40,M40,N40
$ws.Range("H40").Value = 3115.8235
$ws.Range("I40").Value = 2708.577
$ws.Range("J40").Value = 4439.375
$ws.Range("K40").Value = 2708.577
$ws.Range("L40").Value = 4439.375
$ws.Range("M40").Value = -2572.577
$ws.Range("N40").Value = -4711.375
# Row 126: H126,I126,J126,K126,L126,M126,N126
$ws.Range("H126").Value = 5300.222
$ws.Range("I126").Value = 3551
$ws.Range("J126").Value = 6699.6
$ws.Range("K126").Value = 10653
$ws.Range("L126").Value = 20098.8
$ws.Range("M126").Value = -8183
$ws.Range("N126").Value = -25038.8
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 20617.754
$ws.Range("I132").Value = 1519.5
$ws.Range("J132").Value = 93538.37
$ws.Range("K132").Value = 4558.5
$ws.Range("L132").Value = 280615.11
$ws.Range("M132").Value = -2028.5
$ws.Range("N132").Value = -285675.11
# Row 140: H140,J140,L140,N140
$ws.Range("H140").Value = 41151.6
$ws.Range("J140").Value = 41151.6
$ws.Range("L140").Value = 41151.6
$ws.Range("N140").Value = -51511.6
# Row 141: H141,J141,L141,N141
$ws.Range("H141").Value = 57399.645
$ws.Range("J141").Value = 57399.645
$ws.Range("L141").Value = 57399.645
$ws.Range("N141").Value = -67759.64499999999

$ws = $wb.Worksheets.Item("WVR")
# Row 126: H126,I126,J126,K126,L126,M126,N126
$ws.Range("H126").Value = 1100.0476
$ws.Range("I126").Value = 1147.4117
$ws.Range("J126").Value = 898.75
$ws.Range("K126").Value = 3442.2351
$ws.Range("L126").Value = 2696.25
$ws.Range("M126").Value = -972.2351000000003
$ws.Range("N126").Value = -7636.25
